$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure date-like text cells are not auto-converted to Excel date serials
$dateCells = @("AA19", "Y19", "AA20", "Y20", "AA50", "Y50", "AA53", "Y53")
foreach ($addr in $dateCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 18
$ws.Range("A18").Value = 111748673
$ws.Range("B18").Value = 90332
$ws.Range("E18").Value = 4769
$ws.Range("F18").Value = "Svavelriska"
$ws.Range("G18").Value = "Lactarius scrobiculatus"
$ws.Range("H18").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q18").Value = 664574.9395623421
$ws.Range("R18").Value = 6698990.601510798

# Row 19
$ws.Range("A19").Value = 111748762
$ws.Range("B19").Value = 98535
$ws.Range("E19").Value = 222498
$ws.Range("F19").Value = "Blåsippa"
$ws.Range("G19").Value = "Hepatica nobilis"
$ws.Range("H19").Value = "Schreb."
$ws.Range("Q19").Value = 664635.6605944363
$ws.Range("R19").Value = 6698854.271896686
$ws.Range("Y19").Value = "2023-07-01"
$ws.Range("AA19").Value = "2023-07-01"

# Row 20
$ws.Range("A20").Value = 111748702
$ws.Range("B20").Value = 85089
$ws.Range("E20").Value = 3762
$ws.Range("F20").Value = "Olivspindling"
$ws.Range("G20").Value = "Cortinarius venetus"
$ws.Range("H20").Value = "(Fr.:Fr.) Fr."
$ws.Range("Q20").Value = 664852.3607732435
$ws.Range("R20").Value = 6699289.765398038
$ws.Range("Y20").Value = "2023-08-26"
$ws.Range("AA20").Value = "2023-08-26"

# Row 40
$ws.Range("A40").Value = 111748693
$ws.Range("B40").Value = 89183
$ws.Range("E40").Value = 3215
$ws.Range("F40").Value = "Rödgul trumpetsvamp"
$ws.Range("G40").Value = "Craterellus lutescens"
$ws.Range("H40").Value = "(Fr.) Fr."
$ws.Range("Q40").Value = 664610.9204746395
$ws.Range("R40").Value = 6698842.259884536

# Row 41
$ws.Range("A41").Value = 111748771
$ws.Range("B41").Value = 85253
$ws.Range("E41").Value = 3674
$ws.Range("F41").Value = "Anisspindling"
$ws.Range("G41").Value = "Cortinarius odorifer"
$ws.Range("H41").Value = "Britzelm."
$ws.Range("Q41").Value = 664803.0104039316
$ws.Range("R41").Value = 6699372.152713455

# Row 42
$ws.Range("A42").Value = 111748670
$ws.Range("B42").Value = 90332
$ws.Range("E42").Value = 4769
$ws.Range("F42").Value = "Svavelriska"
$ws.Range("G42").Value = "Lactarius scrobiculatus"
$ws.Range("H42").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q42").Value = 664583.6729684926
$ws.Range("R42").Value = 6699322.601555262

# Row 49
$ws.Range("A49").Value = 111748770
$ws.Range("B49").Value = 85253
$ws.Range("E49").Value = 3674
$ws.Range("F49").Value = "Anisspindling"
$ws.Range("G49").Value = "Cortinarius odorifer"
$ws.Range("H49").Value = "Britzelm."
$ws.Range("Q49").Value = 664698.8371583781
$ws.Range("R49").Value = 6699057.103091458

# Row 50
$ws.Range("A50").Value = 111748760
$ws.Range("B50").Value = 98535
$ws.Range("E50").Value = 222498
$ws.Range("F50").Value = "Blåsippa"
$ws.Range("G50").Value = "Hepatica nobilis"
$ws.Range("H50").Value = "Schreb."
$ws.Range("Q50").Value = 664590.5125486635
$ws.Range("R50").Value = 6699346.66743302
$ws.Range("Y50").Value = "2023-07-01"
$ws.Range("AA50").Value = "2023-07-01"

# Row 52
$ws.Range("A52").Value = 111748707
$ws.Range("B52").Value = 96369
$ws.Range("E52").Value = 219862
$ws.Range("F52").Value = "Nästrot"
$ws.Range("G52").Value = "Neottia nidus-avis"
$ws.Range("H52").Value = "(L.) Rich."
$ws.Range("Q52").Value = 664850.52293942
$ws.Range("R52").Value = 6699362.928853855

# Row 53
$ws.Range("A53").Value = 111748695
$ws.Range("B53").Value = 89183
$ws.Range("E53").Value = 3215
$ws.Range("F53").Value = "Rödgul trumpetsvamp"
$ws.Range("G53").Value = "Craterellus lutescens"
$ws.Range("H53").Value = "(Fr.) Fr."
$ws.Range("Q53").Value = 664758.1487707719
$ws.Range("R53").Value = 6698951.896335257
$ws.Range("Y53").Value = "2023-08-26"
$ws.Range("AA53").Value = "2023-08-26"

# Restore default (General) formatting now that text values are locked in as text
foreach ($addr in $dateCells) {
    $ws.Range($addr).ClearFormats()
}